$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row/column header labels: strip " Diff-in-Diff" suffix
$ws.Range("B1").Value = "FFR"
$ws.Range("C1").Value = "C/A"
$ws.Range("D1").Value = "U"
$ws.Range("E1").Value = "`$\pi`$"

$ws.Range("A2").Value = "FFR"
$ws.Range("A3").Value = "C/A"
$ws.Range("A4").Value = "U"
$ws.Range("A5").Value = "`$\pi`$"

# Helper: write a value that must stay text even if it looks numeric.
# A leading apostrophe forces Excel to store it as text, then resetting
# the cell style back to Normal drops the quote-prefix formatting while
# the text stays text.
function Set-TextValue($addr, $text) {
    $r = $ws.Range($addr)
    $r.Value = "'" + $text
    $r.Style = "Normal"
}

# Column B values
Set-TextValue "B3" "0.063"
$ws.Range("B4").Value = "-2.343***"
$ws.Range("B5").Value = "0.519*"
$ws.Range("B6").Value = "-0.363*"
$ws.Range("B7").Value = 0.75

# Column C values
Set-TextValue "C2" "0.408"
Set-TextValue "C4" "-0.989"
Set-TextValue "C5" "0.057"
Set-TextValue "C6" "0.535"
$ws.Range("C7").Value = 0.2

# Column D values
$ws.Range("D2").Value = "-0.296***"
Set-TextValue "D3" "-0.019"
$ws.Range("D5").Value = "0.23**"
Set-TextValue "D6" "0.013"
$ws.Range("D7").Value = 0.77

# Column E values
$ws.Range("E2").Value = "0.301*"
Set-TextValue "E3" "0.005"
$ws.Range("E4").Value = "1.055**"
Set-TextValue "E6" "-0.115"
$ws.Range("E7").Value = 0.17
